# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
# Both sheets carry the same event rows; most F-values move in lockstep,
# except row 22 which ends up one apart between the two sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# row -> new F value, shared by both sheets
$commonUpdates = @{
    2  = 825
    4  = 1147
    5  = 48
    6  = 12354
    7  = 52
    9  = 493
    10 = 442
    11 = 1130
    12 = 911
    13 = 13612
    14 = 13806
    16 = 163
    19 = 1032
    23 = 4903
    24 = 219
}

foreach ($sheet in @($ws1, $ws4)) {
    foreach ($row in $commonUpdates.Keys) {
        $sheet.Cells.Item($row, 6).Value = $commonUpdates[$row]
    }
}

# Row 22 diverges between the two sheets.
$ws1.Cells.Item(22, 6).Value = 132
$ws4.Cells.Item(22, 6).Value = 133
